# TC04_C3DC_phs002371_AnatomicSite-Blood.xlsx
# Update the SQL queries on Sheet1 so the joins use the renamed
# "study_id" / "participant_id" key columns instead of the old
# generic "id" columns, and refresh the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Query([string]$sql) {
    $sql = $sql.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $sql = $sql.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $sql = $sql.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $sql = $sql.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $sql = $sql.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $sql = $sql.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    return $sql
}

# Cells holding SQL text: C2 (StatQuery) plus the TabQuery column (B)
# for every tab row (Studies, Participants, Diagnosis, Treatment,
# TreatmentResp, Survival).
$targets = @(
    @{Row=2; Col=3},
    @{Row=2; Col=2},
    @{Row=3; Col=2},
    @{Row=4; Col=2},
    @{Row=5; Col=2},
    @{Row=6; Col=2},
    @{Row=7; Col=2}
)

foreach ($t in $targets) {
    $cell = $ws.Cells.Item($t.Row, $t.Col)
    $cell.Value = Fix-Query $cell.Text
}

# Reflect the author's final cursor position after the edits.
$ws.Range("C7").Select()
